$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.22594800312889163
$ws.Range("A2").Value = -0.0059999999889441824
$ws.Range("A3").Value = -0.0039999999897428751
$ws.Range("A4").Value = -0.0079999999816564582
$ws.Range("A5").Value = -0.0029999999893730589
$ws.Range("A6").Value = -0.0019999999885804698
$ws.Range("A7").Value = -0.0099999999742639289
$ws.Range("A8").Value = -0.0099999999735231881
$ws.Range("A9").Value = -0.001999999986689982
$ws.Range("A10").Value = -0.0019999999858768547
$ws.Range("A11").Value = -0.0029999999840297775
$ws.Range("A12").Value = -0.0034999999829561368
$ws.Range("A13").Value = 0.0015344482184360686
$ws.Range("A14").Value = -0.0079999999740127947
$ws.Range("A15").Value = -0.00099999998611544072
$ws.Range("A16").Value = -0.0019999999841813221
$ws.Range("A17").Value = -0.0019999999839148686
$ws.Range("A18").Value = -0.0039999999803264075
$ws.Range("A19").Value = -0.0039999999923918672
$ws.Range("A20").Value = -0.0039999999919562157
$ws.Range("A21").Value = -0.0039999999918709506
$ws.Range("A22").Value = -0.003999999991803449
$ws.Range("A23").Value = -0.011533613528118991
$ws.Range("A24").Value = -0.054464233439303733
$ws.Range("A25").Value = -0.019999999958019821
$ws.Range("A26").Value = -0.0024999999874779633
$ws.Range("A27").Value = -0.0024999999871160306
$ws.Range("A28").Value = -0.0019999999865936147
$ws.Range("A29").Value = -0.0069999999768493026
$ws.Range("A30").Value = -0.059999999883495025
$ws.Range("A31").Value = -0.0069999999770953281
$ws.Range("A32").Value = 0.036828604044591984
$ws.Range("A33").Value = -0.0039999999825788279
